$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 824, pushing existing rows 824-882 down to 825-883.
$ws.Rows("824:824").Insert()

# Populate the newly inserted row 824 with the new record.
$ws.Range("A824").Value = 3
$ws.Range("B824").Value = "Femacal de La Calera"
$ws.Range("C824").Value = "Coquimbo"
$ws.Range("D824").Value = 45021
$ws.Range("E824").Value = 5
$ws.Range("F824").Value = 100112045
$ws.Range("G824").Value = "Zapallo"
$ws.Range("H824").Value = "Camote"
$ws.Range("I824").Value = "1a (cosecha)"
$ws.Range("J824").Value = 180
$ws.Range("K824").Value = 500
$ws.Range("L824").Value = 550
$ws.Range("M824").Value = 525
$ws.Range("N824").Value = "$/kilo (volumen en unidades)"
$ws.Range("O824").Value = "Provincia de Talca"
$ws.Range("P824").Value = 525
$ws.Range("Q824").Value = 1
$ws.Range("R824").Value = "Hortaliza"
